# regen sval data to filter save games
# Update the numeric stat columns (TB, d2S, K, IP) and the derived "sum"
# column (B:E and G) for each data row. Column F (Win) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    3  = @{ B = 3.272327238179451;  C = 109.9114832445916;  D = 0.1496068669990043; E = 13.86384647080068;   G = 127.1972638205707 }
    4  = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
    5  = @{ B = 1.445647641019636;  C = 0.3048912486333797; D = 3.223369029078222;  E = 0.5333859586016987;  G = 5.507293877332936 }
    6  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    7  = @{ B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.536033448013082 }
    8  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    9  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    10 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 3.755628166162433 }
    11 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
